# Add "Bye" and "ESPN Projection" columns to the depth chart header row,
# reflecting the new draft-history / depth-chart columns added in this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Bye"
$ws.Range("D1").Value = "ESPN Projection"

# Match the resulting selection/active cell left behind after the edit.
$ws.Range("D1").Select()
